# Predictions_1-26-23.xlsx edit:
#   "Winning Team" columns (in Table2 and Table3) renamed to "Home Result",
#   and their W/L text changed from "Home Wins"/"Away Wins" to "W"/"L".
#   Dependent formula in column M updated to reference the renamed column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the calculated "Winning Team" column headers to "Home Result" ---
# (Table2 header is F1, Table3 header is K1; updating the header cell value
#  renames the underlying table column as well.)
$ws.Range("F1").Value2 = "Home Result"
$ws.Range("K1").Value2 = "Home Result"

# --- Update the Table2 "Home Result" calculated column formula (W / L) ---
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("F$r").Formula = '=IF(Table2[[#This Row],[Home Score]]>Table2[[#This Row],[Away Score]],"W", "L")'
}

# --- Update the Table3 "Home Result" calculated column formula (W / L) ---
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("K$r").Formula = '=IF(Table3[[#This Row],[Actual Home Score]]>Table3[[#This Row],[Actual Away Score]], "W", "L")'
}

# --- Update the "Correct"/"Inccorect" comparison formula to use the renamed column ---
for ($r = 2; $r -le 17; $r++) {
    $ws.Range("M$r").Formula = '=+IF(Table2[[#This Row],[Home Result]]=Table3[[#This Row],[Home Result]], "Correct", "Inccorect")'
}

# --- Update the saved cursor/selection position ---
[void]$ws.Range("I27").Select()
